# Append a new record (row 62) to the Optical_Power sheet, mirroring the
# existing data layout (columns A-N, header in row 1, data starting row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

# Text columns (A-L) -- force Text format so Excel doesn't silently coerce
# numeric-looking strings (e.g. "-477", "807472732", dates) into numbers
# or date serials.
$textRange = "A" + $row + ":L" + $row
$ws.Range($textRange).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "-477"
$ws.Cells.Item($row, 2).Value  = "6/13/2025"
$ws.Cells.Item($row, 3).Value  = "Castañares 4511"
$ws.Cells.Item($row, 4).Value  = "8"
$ws.Cells.Item($row, 5).Value  = "807472732"
$ws.Cells.Item($row, 6).Value  = "Optical Power"
$ws.Cells.Item($row, 7).Value  = "Pendiente"
$ws.Cells.Item($row, 8).Value  = ""
$ws.Cells.Item($row, 9).Value  = "0"
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Terminal"

# Numeric coordinate columns (M, N)
$ws.Cells.Item($row, 13).Value = -58.470376
$ws.Cells.Item($row, 14).Value = -34.664751
